# Scheduled runner update: refresh market-price-derived profit figures
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ and LeveProfit NQ/HQ
# columns H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR
# sheets with newly pulled marketboard data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 861.5454999999999
$ws.Range("I6").Value = 954.7778
$ws.Range("K6").Value = 2864.3334
$ws.Range("M6").Value = -2752.3334

$ws.Range("H29").Value = 8575
$ws.Range("J29").Value = 8575
$ws.Range("L29").Value = 25725
$ws.Range("N29").Value = -26287

$ws.Range("H38").Value = 1505.1666
$ws.Range("I38").Value = 1206.2
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 3618.6
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = -3246.6
$ws.Range("N38").Value = -9744

$ws.Range("H58").Value = 1203
$ws.Range("J58").Value = 3796
$ws.Range("L58").Value = 11388
$ws.Range("N58").Value = -11688

$ws.Range("H61").Value = 1887.4286
$ws.Range("I61").Value = 1887.4286
$ws.Range("K61").Value = 5662.2858
$ws.Range("M61").Value = -5490.2858

$ws.Range("H87").Value = 49900
$ws.Range("J87").Value = 49900
$ws.Range("L87").Value = 49900
$ws.Range("N87").Value = -52396

$ws.Range("H88").Value = 1332.5
$ws.Range("I88").Value = 1665
$ws.Range("J88").Value = 1000
$ws.Range("K88").Value = 1665
$ws.Range("L88").Value = 1000
$ws.Range("M88").Value = -1259
$ws.Range("N88").Value = -1812

$ws.Range("H90").Value = 49900
$ws.Range("J90").Value = 49900
$ws.Range("L90").Value = 149700
$ws.Range("N90").Value = -162180

$ws.Range("H91").Value = 1332.5
$ws.Range("I91").Value = 1665
$ws.Range("J91").Value = 1000
$ws.Range("K91").Value = 1665
$ws.Range("L91").Value = 1000
$ws.Range("M91").Value = -261
$ws.Range("N91").Value = -3808

$ws.Range("H116").Value = 12625.25
$ws.Range("I116").Value = 22600.6
$ws.Range("J116").Value = 5500
$ws.Range("K116").Value = 22600.6
$ws.Range("L116").Value = 5500
$ws.Range("M116").Value = -19158.6
$ws.Range("N116").Value = -12384

$ws.Range("H137").Value = 56687.168
$ws.Range("I137").Value = 744.3
$ws.Range("J137").Value = 126615.75
$ws.Range("K137").Value = 2232.9
$ws.Range("L137").Value = 379847.25
$ws.Range("M137").Value = 317.1000000000004
$ws.Range("N137").Value = -384947.25

$ws.Range("H138").Value = 3195.83
$ws.Range("I138").Value = 3211
$ws.Range("J138").Value = 3183.276
$ws.Range("K138").Value = 9633
$ws.Range("L138").Value = 9549.828
$ws.Range("M138").Value = -4493
$ws.Range("N138").Value = -19829.828

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2431.0852
$ws.Range("I32").Value = 1912.25
$ws.Range("K32").Value = 1912.25
$ws.Range("M32").Value = -1625.25

$ws.Range("H45").Value = 18001700
$ws.Range("I45").Value = 22501500
$ws.Range("K45").Value = 22501500
$ws.Range("M45").Value = -22501123

$ws.Range("H61").Value = 4292.636
$ws.Range("I61").Value = 1043.3334
$ws.Range("K61").Value = 1043.3334
$ws.Range("M61").Value = -831.3334

$ws.Range("H74").Value = 842.46155
$ws.Range("I74").Value = 761.087
$ws.Range("J74").Value = 1466.3334
$ws.Range("K74").Value = 761.087
$ws.Range("L74").Value = 1466.3334
$ws.Range("M74").Value = 112.913
$ws.Range("N74").Value = -3214.3334

$ws.Range("H77").Value = 842.46155
$ws.Range("I77").Value = 761.087
$ws.Range("J77").Value = 1466.3334
$ws.Range("K77").Value = 3805.435
$ws.Range("L77").Value = 7331.666999999999
$ws.Range("M77").Value = 562.5650000000001
$ws.Range("N77").Value = -16067.667

$ws.Range("H132").Value = 2357.4324
$ws.Range("I132").Value = 2041.8572
$ws.Range("J132").Value = 2771.625
$ws.Range("K132").Value = 6125.571599999999
$ws.Range("L132").Value = 8314.875
$ws.Range("M132").Value = -3595.571599999999
$ws.Range("N132").Value = -13374.875

$ws.Range("H136").Value = 4292.636
$ws.Range("I136").Value = 1043.3334
$ws.Range("K136").Value = 3130.0002
$ws.Range("M136").Value = -580.0001999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 46952
$ws.Range("J76").Value = 46952
$ws.Range("L76").Value = 46952
$ws.Range("N76").Value = -47582

$ws.Range("H79").Value = 46952
$ws.Range("J79").Value = 46952
$ws.Range("L79").Value = 46952
$ws.Range("N79").Value = -49136

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9583.333000000001
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = -4888
$ws.Range("N4").Value = -10224

$ws.Range("H31").Value = 1985.0638
$ws.Range("I31").Value = 1112.35
$ws.Range("J31").Value = 2631.5186
$ws.Range("K31").Value = 1112.35
$ws.Range("L31").Value = 2631.5186
$ws.Range("M31").Value = -817.3499999999999
$ws.Range("N31").Value = -3221.5186

$ws.Range("H34").Value = 1985.0638
$ws.Range("I34").Value = 1112.35
$ws.Range("J34").Value = 2631.5186
$ws.Range("K34").Value = 1112.35
$ws.Range("L34").Value = 2631.5186
$ws.Range("M34").Value = -910.3499999999999
$ws.Range("N34").Value = -3035.5186

$ws.Range("H58").Value = 1893434.1
$ws.Range("I58").Value = 3346229.8
$ws.Range("J58").Value = 4799.7
$ws.Range("K58").Value = 3346229.8
$ws.Range("L58").Value = 4799.7
$ws.Range("M58").Value = -3346026.8
$ws.Range("N58").Value = -5205.7

$ws.Range("H99").Value = 1430656.6
$ws.Range("I99").Value = 2501574
$ws.Range("K99").Value = 2501574
$ws.Range("M99").Value = -2500076

$ws.Range("H105").Value = 1031.2222
$ws.Range("I105").Value = 910.25
$ws.Range("K105").Value = 910.25
$ws.Range("M105").Value = 836.75

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H126").Value = 1430656.6
$ws.Range("I126").Value = 2501574
$ws.Range("K126").Value = 7504722
$ws.Range("M126").Value = -7502252

$ws.Range("H132").Value = 4925.35
$ws.Range("I132").Value = 3812.25
$ws.Range("J132").Value = 5667.4165
$ws.Range("K132").Value = 11436.75
$ws.Range("L132").Value = 17002.2495
$ws.Range("M132").Value = -8906.75
$ws.Range("N132").Value = -22062.2495

$ws.Range("H136").Value = 1893434.1
$ws.Range("I136").Value = 3346229.8
$ws.Range("J136").Value = 4799.7
$ws.Range("K136").Value = 10038689.4
$ws.Range("L136").Value = 14399.1
$ws.Range("M136").Value = -10036139.4
$ws.Range("N136").Value = -19499.1

$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1460
$ws.Range("J132").Value = 1850
$ws.Range("L132").Value = 16650
$ws.Range("N132").Value = -21710

$ws.Range("H133").Value = 83335336
$ws.Range("I133").Value = 83335336
$ws.Range("K133").Value = 250006008
$ws.Range("M133").Value = -250000948

$ws.Range("H137").Value = 2903
$ws.Range("J137").Value = 4375
$ws.Range("L137").Value = 13125
$ws.Range("N137").Value = -23325

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H132").Value = 803393.2
$ws.Range("I132").Value = 1132677
$ws.Range("J132").Value = 3703.8572
$ws.Range("K132").Value = 3398031
$ws.Range("L132").Value = 11111.5716
$ws.Range("M132").Value = -3395501
$ws.Range("N132").Value = -16171.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5218.364
$ws.Range("I7").Value = 2486.1428
$ws.Range("K7").Value = 2486.1428
$ws.Range("M7").Value = -2374.1428

$ws.Range("H46").Value = 1077.762
$ws.Range("I46").Value = 489
$ws.Range("J46").Value = 1372.1428
$ws.Range("K46").Value = 489
$ws.Range("L46").Value = 1372.1428
$ws.Range("M46").Value = -301
$ws.Range("N46").Value = -1748.1428

$ws.Range("H61").Value = 3300
$ws.Range("I61").Value = 2450
$ws.Range("K61").Value = 2450
$ws.Range("M61").Value = -2248

$ws.Range("H113").Value = 3300
$ws.Range("I113").Value = 2450
$ws.Range("K113").Value = 2450
$ws.Range("M113").Value = -280

$ws.Range("H122").Value = 5220.115
$ws.Range("I122").Value = 4371.077
$ws.Range("J122").Value = 6069.154
$ws.Range("K122").Value = 13113.231
$ws.Range("L122").Value = 18207.462
$ws.Range("M122").Value = -10663.231
$ws.Range("N122").Value = -23107.462

$ws.Range("H126").Value = 5218.364
$ws.Range("I126").Value = 2486.1428
$ws.Range("K126").Value = 7458.428400000001
$ws.Range("M126").Value = -4988.428400000001

$ws.Range("H136").Value = 6659
$ws.Range("I136").Value = 4475
$ws.Range("J136").Value = 8115
$ws.Range("K136").Value = 13425
$ws.Range("L136").Value = 24345
$ws.Range("M136").Value = -10875
$ws.Range("N136").Value = -29445

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 17000
$ws.Range("J31").Value = 17000
$ws.Range("L31").Value = 17000
$ws.Range("N31").Value = -17696

$ws.Range("H126").Value = 3310.0881
$ws.Range("I126").Value = 2669.25
$ws.Range("K126").Value = 8007.75
$ws.Range("M126").Value = -5537.75

$ws.Range("H132").Value = 1722.6765
$ws.Range("I132").Value = 736.5714
$ws.Range("J132").Value = 2412.95
$ws.Range("K132").Value = 2209.7142
$ws.Range("L132").Value = 7238.849999999999
$ws.Range("M132").Value = 320.2857999999997
$ws.Range("N132").Value = -12298.85

$ws.Range("H136").Value = 13230073
$ws.Range("I136").Value = 34725076
$ws.Range("J136").Value = 2377.9614
$ws.Range("K136").Value = 104175228
$ws.Range("L136").Value = 7133.8842
$ws.Range("M136").Value = -104172678
$ws.Range("N136").Value = -12233.8842
